$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (AC1) into the three
# new header cells so they pick up the same bold/centered/bordered style,
# then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2..53 get the team record values (Wins, Losses, Ties)
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 94  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 68  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
